$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 395; this shifts rows 395:464 down to 396:465
# (matching the dimension change A1:T464 -> A1:T465 in the diff).
$ws.Rows("395:395").Insert()

# Populate the newly inserted row 395 with the new record's data.
$ws.Range("A395").Value = 9
$ws.Range("B395").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C395").Value = "Metropolitana"
$ws.Range("D395").Value = 44776
$ws.Range("E395").Value = 13
$ws.Range("F395").Value = "Fruta"
$ws.Range("G395").Value = 100108
$ws.Range("H395").Value = "Tropicales y subtropicales"
$ws.Range("I395").Value = 100108002
$ws.Range("J395").Value = "Mango"
$ws.Range("K395").Value = "Sin especificar"
$ws.Range("L395").Value = "Primera"
$ws.Range("M395").Value = 630
$ws.Range("N395").Value = 8500
$ws.Range("O395").Value = 9500
$ws.Range("P395").Value = 8944
$ws.Range("Q395").Value = "$/bandeja 4 kilos"
$ws.Range("R395").Value = "México"
$ws.Range("S395").Value = 2236
$ws.Range("T395").Value = 4
